# Auto-generated edit script updating the cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.750.39'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '3.452.25'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.87'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.68'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.618'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +13.05%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '3.451.84'
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.23'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.447'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.43%  '
$ws.Range('D13').Value = '4.053.57'
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000192'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('E16').Value = '  +3.32%  '
$ws.Range('D17').Value = '64.823.34'
$ws.Range('E17').Value = '  +1.40%  '
$ws.Range('D18').Value = '3.466.74'
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.45'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.32'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.20'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.12'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.551'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.65%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.58'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.95%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000119'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.93'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +5.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.178'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.51'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +10.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.12'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.61'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.17'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.61'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +10.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '161.40'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('E37').Value = '  +4.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0775'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.28%  '
$ws.Range('D39').Value = '2.953.12'
$ws.Range('E39').Value = '  -1.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.56'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.48%  '
$ws.Range('E41').Value = '  +6.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.61'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.23%  '
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('E44').Value = '  +1.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.772'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.70'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +10.51%  '
$ws.Range('E47').Value = '  +1.54%  '
$ws.Range('E48').Value = '  +8.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '309.60'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.38%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.867'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.65%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.61'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.28%  '
